$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 5 values to the new (rounded) figures
$ws.Range("B5").Value = 7.55
$ws.Range("C5").Value = 5.82
$ws.Range("D5").Value = 0.81
$ws.Range("E5").Value = 16.39
$ws.Range("F5").Value = 13.4
$ws.Range("G5").Value = 6.23
$ws.Range("H5").Value = 28.47
$ws.Range("I5").Value = 9.31
$ws.Range("J5").Value = 4.05
$ws.Range("K5").Value = 6.05
$ws.Range("L5").Value = 6.61
$ws.Range("M5").Value = 6.88
$ws.Range("N5").Value = 1.87
$ws.Range("O5").Value = 6.02
$ws.Range("P5").Value = 8.53
$ws.Range("Q5").Value = 5.19
$ws.Range("R5").Value = 0.76
$ws.Range("S5").Value = 0.42
$ws.Range("T5").Value = 84.54
$ws.Range("U5").Value = 17.05
$ws.Range("V5").Value = 5.55
$ws.Range("W5").Value = 11.33
$ws.Range("X5").Value = 6.07
$ws.Range("Y5").Value = 0.77
$ws.Range("Z5").Value = 13.12
$ws.Range("AA5").Value = 4.88
$ws.Range("AB5").Value = 4.36
$ws.Range("AC5").Value = 5.04
$ws.Range("AD5").Value = 7
$ws.Range("AE5").Value = 0.55
$ws.Range("AF5").Value = 26.13
$ws.Range("AG5").Value = 3.09
$ws.Range("AH5").Value = 6.92

# Remove row 6 entirely (dataset now ends at row 5)
$ws.Rows("6:6").Delete()
